$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The subtitle placeholder on slide 1 moves down/shrinks and gains a line of
# text that points people at the GitHub repo with the examples.
$subtitle = $s.Shapes.Item(2)

$subtitle.Left = 120
$subtitle.Top = 438.031
$subtitle.Width = 720
$subtitle.Height = 60.78732283464567

$tr = $subtitle.TextFrame.TextRange
$tr.Text = "Most of the examples here are available at https://github.com/bassmandja/docker101"

# Turn just the URL portion of the new text into a hyperlink.
$urlStart = ("Most of the examples here are available at ").Length + 1
$urlLength = ("https://github.com/bassmandja/docker101").Length
$linkRange = $tr.Characters($urlStart, $urlLength)
$linkRange.ActionSettings.Item(1).Hyperlink.Address = "https://github.com/bassmandja/docker101"
